$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0) The document has a hidden "_GoBack" bookmark sitting at the end of the
#    "There were no bugs found..." paragraph. The edit relocates it to the
#    end of the paragraph touched in step 1, so remove the old one now
#    (bookmark names must stay unique) and re-add it in its new home below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 1) Rewrite the "In the second deliverable, ..." sentence through the end of
#    that paragraph, splitting it into the many small runs that the author's
#    edit produced, adding a proofErr spell-check wrap around "checkstyle",
#    and moving the _GoBack bookmark to the end of this paragraph.
# ---------------------------------------------------------------------------
$findStart = $d.Content.Duplicate
$findStart.Find.Execute("In the second deliverable, we started ") | Out-Null

$findEnd = $d.Content.Duplicate
$findEnd.Find.Execute("so we would test the expression counter before we would test the Halstead metrics. ") | Out-Null

$target = $d.Range($findStart.Start, $findEnd.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">In the second deliverable, we </w:t></w:r><w:r><w:t xml:space="preserve">first </w:t></w:r><w:r><w:t xml:space="preserve">started to </w:t></w:r><w:r><w:t xml:space="preserve">test the overall code and </w:t></w:r><w:r><w:t xml:space="preserve">then </w:t></w:r><w:r><w:t>started testing individual component</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve">. The order that we tested </w:t></w:r><w:r><w:t>our components</w:t></w:r><w:r><w:t xml:space="preserve"> depended </w:t></w:r><w:r><w:t xml:space="preserve">on if there were dependencies. We first tested the variable counter, comment counter, cast counter, expression counter, looping counter, and the method counter because they did not depend on any other classes expect for the eclipse </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>checkstyle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class. This would be first level in the call graph. The second level of the call graph would be the Halsted Metrics class because it depends on the expression counter. The third level of the call graph would be the Maintainability Index class because it depends on the comment counter and the Halsted Metrics. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) Merge the duplicated "Bugs discovered during un" + " " +
#    "Bugs discovered during unit testing: " runs into a single run (same
#    visible text, just collapsed back into one run).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Bugs discovered during un Bugs discovered during unit testing: ", $true, $false, $false, $false, $false, $true, 1, $false, "Bugs discovered during un Bugs discovered during unit testing: ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Add a lastRenderedPageBreak marker before "While testing the " (keeping
#    the run's existing character formatting).
# ---------------------------------------------------------------------------
$wt = $d.Content.Duplicate
$wt.Find.Execute("While testing the ") | Out-Null
$wtRange = $d.Range($wt.Start, $wt.End)
$wtXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">While testing the </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$wtRange.InsertXML($wtXml)

